$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force text storage for values that look like numbers/dates
# by setting the NumberFormat to "@" (text) before assignment, then resetting
# the cell style back to "Normal" so no stray style index is left on the cell.

# Row 2
$ws.Range("A2").Value = 1222

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2000/12/23"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "Groceries"
$ws.Range("D2").Value = "lmao"
$ws.Range("E2").Value = "Cash"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2024/06/01"
$ws.Range("F2").Style = "Normal"

# Row 3
$ws.Range("A3").Value = 121212

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2000/12/22"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = "Groceries"
$ws.Range("D3").Value = "sdaffafaasdf"
$ws.Range("E3").Value = "Cash"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2024/06/01"
$ws.Range("F3").Style = "Normal"

# Row 4
$ws.Range("A4").Value = 12121221

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2000/12/22"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").Value = "Groceries"
$ws.Range("D4").Value = "erfan"
$ws.Range("E4").Value = "Cash"

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "2024/05/27"
$ws.Range("F4").Style = "Normal"

# Row 5 - A5 is a string value "23111" in the target, not numeric
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "23111"
$ws.Range("A5").Style = "Normal"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2000/12/22"
$ws.Range("B5").Style = "Normal"

$ws.Range("C5").Value = "Groceries"
$ws.Range("D5").Value = "asdffdassdf"
$ws.Range("E5").Value = "Cash"

$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "2024/05/30"
$ws.Range("F5").Style = "Normal"
